$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H) — reuse the existing header style (same as
# G1's "sum" header) by copying formats only, then set the values.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
